# Weekly crypto-price/volume refresh (GitHub Actions bot).
# D = Price, E = Volume(1h); both columns are plain text cells in the sheet
# (prices use "." as a thousands separator in some rows, e.g. "66.639.46"),
# so numeric-looking values are written with a leading "`'" to force Excel
# to keep them as text instead of parsing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2: Bitcoin
$ws.Range("D2").Value = "66.639.46"
$ws.Range("E2").Value = "  -0.19%  "

# row 3: Ethereum
$ws.Range("D3").Value = "3.808.12"
$ws.Range("E3").Value = "  +0.19%  "

# row 4: TetherUSD
$ws.Range("D4").Value = "`'0.998"
$ws.Range("E4").Value = "  -0.08%  "

# row 5: BNB
$ws.Range("D5").Value = "`'435.19"
$ws.Range("E5").Value = "  +5.38%  "

# row 6: Solana
$ws.Range("D6").Value = "`'139.10"
$ws.Range("E6").Value = "  +4.65%  "

# row 7: XRP
$ws.Range("D7").Value = "`'0.627"
$ws.Range("E7").Value = "  +1.91%  "

# row 8: USDC
$ws.Range("E8").Value = "  -0.10%  "

# row 9: Cardano
$ws.Range("D9").Value = "`'0.741"
$ws.Range("E9").Value = "  -0.25%  "

# row 10: Dogecoin
$ws.Range("D10").Value = "`'0.154"
$ws.Range("E10").Value = "  -8.99%  "

# row 11: ShibaInu
$ws.Range("D11").Value = "`'0.0000324"
$ws.Range("E11").Value = "  -12.75%  "

# row 12: Avalanche
$ws.Range("D12").Value = "`'42.88"
$ws.Range("E12").Value = "  +4.18%  "

# row 13: Polkadot
$ws.Range("D13").Value = "`'10.47"
$ws.Range("E13").Value = "  +2.29%  "

# row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.392.79"
$ws.Range("E14").Value = "  +0.24%  "

# row 15: Uniswap
$ws.Range("D15").Value = "`'15.09"
$ws.Range("E15").Value = "  +1.73%  "

# row 16: TRON
$ws.Range("E16").Value = "  -0.13%  "

# row 17: WrappedEther
$ws.Range("D17").Value = "3.822.77"
$ws.Range("E17").Value = "  +0.58%  "

# row 18: Chainlink
$ws.Range("D18").Value = "`'20.02"
$ws.Range("E18").Value = "  +2.24%  "

# row 19: Polygon
$ws.Range("E19").Value = "  +5.19%  "

# row 20: WrappedBTC
$ws.Range("D20").Value = "66.620.94"
$ws.Range("E20").Value = "  -0.59%  "

# row 21: BitcoinCash
$ws.Range("D21").Value = "`'411.57"
$ws.Range("E21").Value = "  -0.68%  "

# row 22: InternetComputer(DFINITY)
$ws.Range("D22").Value = "`'14.85"
$ws.Range("E22").Value = "  -0.63%  "

# row 23: ImmutableX
$ws.Range("D23").Value = "`'3.25"
$ws.Range("E23").Value = "  +5.28%  "

# row 24: Litecoin
$ws.Range("D24").Value = "`'85.29"
$ws.Range("E24").Value = "  -0.89%  "

# row 25: EthereumClassic
$ws.Range("D25").Value = "`'37.17"
$ws.Range("E25").Value = "  +2.08%  "

# row 26: PancakeSwap
$ws.Range("D26").Value = "`'3.36"
$ws.Range("E26").Value = "  +6.67%  "

# row 27: RenderToken
$ws.Range("D27").Value = "`'9.76"
$ws.Range("E27").Value = "  +31.69%  "

# row 28: LEO
$ws.Range("D28").Value = "`'5.57"
$ws.Range("E28").Value = "  -2.57%  "

# row 29: Filecoin
$ws.Range("D29").Value = "`'9.91"
$ws.Range("E29").Value = "  +4.63%  "

# row 30: Hedera
$ws.Range("D30").Value = "`'0.137"
$ws.Range("E30").Value = "  +11.36%  "

# row 31: Cosmos
$ws.Range("D31").Value = "`'13.88"
$ws.Range("E31").Value = "  +11.19%  "

# row 32: Bittensor
$ws.Range("D32").Value = "`'710.10"
$ws.Range("E32").Value = "  +1.96%  "

# row 33: Toncoin
$ws.Range("D33").Value = "`'2.76"
$ws.Range("E33").Value = "  +0.35%  "

# row 34: InjectiveProtocol
$ws.Range("D34").Value = "`'41.79"
$ws.Range("E34").Value = "  +6.38%  "

# row 35: Dai
$ws.Range("E35").Value = "  +0.02%  "

# row 36: Kaspa
$ws.Range("D36").Value = "`'0.150"
$ws.Range("E36").Value = "  -3.66%  "

# row 37 becomes NEARProtocol
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "`'5.61"
$ws.Range("E37").Value = "  +29.99%  "

# row 38 becomes OKB
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "`'56.15"
$ws.Range("E38").Value = "  +1.17%  "

# row 39: VeChain
$ws.Range("D39").Value = "`'0.0475"
$ws.Range("E39").Value = "  +2.51%  "

# row 40: Fetch.AI
$ws.Range("D40").Value = "`'2.77"
$ws.Range("E40").Value = "  +37.61%  "

# row 41: ThetaToken
$ws.Range("D41").Value = "`'2.91"
$ws.Range("E41").Value = "  +0.32%  "

# row 42 becomes PEPE
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0698"
$ws.Range("E42").Value = "  -9.77%  "

# row 43 becomes Stellar
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "`'0.143"
$ws.Range("E43").Value = "  +3.92%  "

# row 44 becomes FirstDigitalUSD
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "`'0.999"
$ws.Range("E44").Value = "  +0.26%  "

# row 45: TheGraph
$ws.Range("D45").Value = "`'0.327"
$ws.Range("E45").Value = "  +9.60%  "

# row 46: ApeXProtocol
$ws.Range("D46").Value = "`'3.24"
$ws.Range("E46").Value = "  +2.14%  "

# row 47: LidoDAOToken
$ws.Range("E47").Value = "  +1.18%  "

# row 48: WEMIXToken
$ws.Range("D48").Value = "`'2.69"
$ws.Range("E48").Value = "  +3.88%  "

# row 49: ARBITRUM
$ws.Range("D49").Value = "`'2.08"
$ws.Range("E49").Value = "  -0.74%  "

# row 50: Monero
$ws.Range("D50").Value = "`'140.88"
$ws.Range("E50").Value = "  -4.37%  "

# row 51: Stacks
$ws.Range("D51").Value = "`'2.81"
$ws.Range("E51").Value = "  -2.74%  "
